# Apply updated values to the "Crédito disponível" worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K9").Value = 18406.99

$ws.Range("N11").Value = 68001.31

$ws.Range("N12").Value = 868.02

$ws.Range("N15").Value = 2058.54
$ws.Range("O15").Value = 1917.94

$ws.Range("N17").Value = 566485.69
$ws.Range("O17").Value = 518528.56

$ws.Range("K18").Value = 0

$ws.Range("O19").Value = 16089.92

$ws.Range("N23").Value = 183561.09
$ws.Range("O23").Value = 165174.7
